$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('C2').Value = 24
$ws.Range('F2').Value = 31
$ws.Range('H2').Value = 46
$ws.Range('K2').Value = 62
$ws.Range('C3').Value = 41
$ws.Range('I3').Value = 80
$ws.Range('D3').Value = 59
$ws.Range('H3').Value = 38
$ws.Range('K3').Value = 88
$ws.Range('F3').Value = 57
$ws.Range('C9').Value = 204
$ws.Range('D9').Value = 179
$ws.Range('I9').Value = 233
$ws.Range('G9').Value = 210
$ws.Range('F9').Value = 212
$ws.Range('H9').Value = 163
$ws.Range('K9').Value = 195
$ws.Range('J9').Value = 174
$ws.Range('B9').Value = 161
$ws.Range('J10').Value = 305
$ws.Range('H10').Value = 186
$ws.Range('K10').Value = 291
$ws.Range('C10').Value = 546
$ws.Range('F10').Value = 946
$ws.Range('I10').Value = 340
$ws.Range('B10').Value = 452
$ws.Range('E10').Value = 813
$ws.Range('D10').Value = 731
$ws.Range('I11').Value = 703
$ws.Range('J11').Value = 626
$ws.Range('F11').Value = 1249
$ws.Range('H11').Value = 443
$ws.Range('E11').Value = 1093
$ws.Range('K11').Value = 652
$ws.Range('G11').Value = 832
$ws.Range('C11').Value = 821
$ws.Range('B11').Value = 671
$ws.Range('D11').Value = 1010

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('D3').Value = 1
$ws.Range('C6').Value = 13
$ws.Range('F6').Value = 24
$ws.Range('G6').Value = 25
$ws.Range('I6').Value = 15
$ws.Range('E7').Value = 29
$ws.Range('E8').Value = 48
$ws.Range('D8').Value = 42
$ws.Range('G8').Value = 57
$ws.Range('I8').Value = 39
$ws.Range('C8').Value = 36
$ws.Range('F8').Value = 78

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('C7').Value = 20
$ws.Range('K7').Value = 12
$ws.Range('C8').Value = 26
$ws.Range('K8').Value = 18

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('F2').Value = 4
$ws.Range('H7').Value = 10
$ws.Range('H8').Value = 20
$ws.Range('F8').Value = 33

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('F3').Value = 6
$ws.Range('F9').Value = 286
$ws.Range('C9').Value = 98
$ws.Range('E9').Value = 239
$ws.Range('B9').Value = 57
$ws.Range('C10').Value = 122
$ws.Range('E10').Value = 282
$ws.Range('F10').Value = 320
$ws.Range('B10').Value = 74

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('K5').Value = 4
$ws.Range('J6').Value = 3
$ws.Range('H6').Value = 3
$ws.Range('K7').Value = 8
$ws.Range('J7').Value = 12
$ws.Range('H7').Value = 4

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range('F6').Value = 51
$ws.Range('F7').Value = 58

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('I6').Value = 10
$ws.Range('I8').Value = 18

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('H5').Value = 4
$ws.Range('J5').Value = 12
$ws.Range('K5').Value = 8
$ws.Range('C7').Value = 17
$ws.Range('I7').Value = 6
$ws.Range('B8').Value = 35
$ws.Range('F8').Value = 50
$ws.Range('J8').Value = 32
$ws.Range('E16').Value = 2
$ws.Range('C19').Value = 26
$ws.Range('K19').Value = 18
$ws.Range('K28').Value = 37
$ws.Range('K29').Value = 8
$ws.Range('C30').Value = 10
$ws.Range('C32').Value = 36
$ws.Range('I32').Value = 39
$ws.Range('D32').Value = 42
$ws.Range('F32').Value = 78
$ws.Range('G32').Value = 57
$ws.Range('E32').Value = 48
$ws.Range('B35').Value = 11
$ws.Range('H36').Value = 20
$ws.Range('F36').Value = 33
$ws.Range('H41').Value = 4
$ws.Range('J42').Value = 8
$ws.Range('D43').Value = 10
$ws.Range('D47').Value = 22
$ws.Range('H47').Value = 18
$ws.Range('B53').Value = 74
$ws.Range('F53').Value = 320
$ws.Range('E53').Value = 282
$ws.Range('C53').Value = 122
$ws.Range('I65').Value = 18
$ws.Range('C67').Value = 5
$ws.Range('F70').Value = 58
$ws.Range('D74').Value = 33
$ws.Range('K76').Value = 24
$ws.Range('F77').Value = 26
$ws.Range('I85').Value = 3
$ws.Range('F86').Value = 13
$ws.Range('K86').Value = 15
$ws.Range('C87').Value = 7
$ws.Range('F88').Value = 5
$ws.Range('H91').Value = 6
$ws.Range('E94').Value = 38
$ws.Range('J94').Value = 7
$ws.Range('D96').Value = 12
$ws.Range('K96').Value = 6
$ws.Range('I97').Value = 3
$ws.Range('J98').Value = 626
$ws.Range('B98').Value = 671
$ws.Range('D98').Value = 1010
$ws.Range('F98').Value = 1249
$ws.Range('I98').Value = 703
$ws.Range('G98').Value = 832
$ws.Range('H98').Value = 443
$ws.Range('E98').Value = 1093
$ws.Range('K98').Value = 652
$ws.Range('C98').Value = 821

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('F5').Value = 5
$ws.Range('F6').Value = 5

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('H3').Value = 1
$ws.Range('H7').Value = 4
$ws.Range('K2').Value = 3

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('F6').Value = 8
$ws.Range('F7').Value = 13
$ws.Range('K7').Value = 15

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K3').Value = 9
$ws.Range('K8').Value = 37

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('H6').Value = 8
$ws.Range('D7').Value = 18
$ws.Range('D8').Value = 22
$ws.Range('H8').Value = 18
$ws.Range('K6').Value = 3

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('K7').Value = 8

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('K6').Value = 5
$ws.Range('K8').Value = 24

$ws = $wb.Worksheets.Item('River North')
$ws.Range('D5').Value = 28
$ws.Range('D6').Value = 33

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('I5').Value = 2
$ws.Range('I6').Value = 3

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('H2').Value = 3
$ws.Range('H9').Value = 6

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range('C6').Value = 5
$ws.Range('C7').Value = 5

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('B5').Value = 8
$ws.Range('B6').Value = 11

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('K3').Value = 2
$ws.Range('D6').Value = 7
$ws.Range('D7').Value = 12
$ws.Range('K7').Value = 6

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('F8').Value = 20
$ws.Range('F9').Value = 26

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('J5').Value = 2
$ws.Range('E6').Value = 35
$ws.Range('E7').Value = 38
$ws.Range('J7').Value = 7

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('C2').Value = 1
$ws.Range('I3').Value = 2
$ws.Range('C7').Value = 17
$ws.Range('I7').Value = 6

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range('I5').Value = 1
$ws.Range('I7').Value = 3

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('J5').Value = 5
$ws.Range('J6').Value = 8

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('C5').Value = 8
$ws.Range('C6').Value = 10

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range('D5').Value = 2
$ws.Range('D6').Value = 2

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('B6').Value = 10
$ws.Range('F6').Value = 11
$ws.Range('J6').Value = 11
$ws.Range('B8').Value = 35
$ws.Range('F8').Value = 50
$ws.Range('J8').Value = 32

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('D5').Value = 3
$ws.Range('D7').Value = 10

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('C3').Value = 1
$ws.Range('C7').Value = 7
